$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Roll the price-list date (A1) forward one month
$ws.Range("A1").Value = 45436

# Update unit prices in column D (rows 31-38)
$ws.Range("D31").Value = 9938.521000000001
$ws.Range("D32").Value = 8865.171
$ws.Range("D33").Value = 8453.887000000001
$ws.Range("D34").Value = 14397.439
$ws.Range("D35").Value = 10743.533
$ws.Range("D36").Value = 10126.608
$ws.Range("D37").Value = 9489.620000000001
$ws.Range("D38").Value = 15004.331
